$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.903.37"
$ws.Range("E2").Value = "  +0.43%  "
$ws.Range("D3").Value = "1.552.23"
$ws.Range("E3").Value = "  +0.45%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.32%  "
$ws.Range("D5").Value = "'206.21"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.88%  "
$ws.Range("D6").Value = "'0.483"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.46%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.34%  "
$ws.Range("B8").Value = "Solana"
$ws.Range("C8").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D8").Value = "'21.53"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.00%  "
$ws.Range("B9").Value = "Cardano"
$ws.Range("C9").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D9").Value = "'0.247"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.80%  "
$ws.Range("E10").Value = "  +0.41%  "
$ws.Range("D11").Value = "'0.0858"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.25%  "
$ws.Range("D12").Value = "1.772.70"
$ws.Range("E12").Value = "  +0.36%  "
$ws.Range("D13").Value = "1.524.73"
$ws.Range("E13").Value = "  -1.76%  "
$ws.Range("D15").Value = "'0.514"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.08%  "
$ws.Range("D16").Value = "26.900.90"
$ws.Range("E16").Value = "  +0.45%  "
$ws.Range("D17").Value = "'61.62"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.18%  "
$ws.Range("D18").Value = "'213.89"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.20%  "
$ws.Range("E19").Value = "  +0.94%  "
$ws.Range("D20").Value = "'7.24"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.10%  "
$ws.Range("E21").Value = "  +0.38%  "
$ws.Range("E22").Value = "  -0.73%  "
$ws.Range("E23").Value = "  +1.56%  "
$ws.Range("D24").Value = "'1.97"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.70%  "
$ws.Range("D25").Value = "'153.30"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.65%  "
$ws.Range("D26").Value = "'6.67"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.74%  "
$ws.Range("D27").Value = "'14.87"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.56%  "
$ws.Range("E28").Value = "  +0.33%  "
$ws.Range("E29").Value = "  +1.74%  "
$ws.Range("E30").Value = "  -0.32%  "
$ws.Range("E31").Value = "  -0.76%  "
$ws.Range("D32").Value = "'3.23"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.01%  "
$ws.Range("D33").Value = "1.369.65"
$ws.Range("E33").Value = "  +0.88%  "
$ws.Range("D34").Value = "'2.96"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.03%  "
$ws.Range("D35").Value = "'1.55"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.89%  "
$ws.Range("D36").Value = "'0.974"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +6.83%  "
$ws.Range("E37").Value = "  +0.57%  "
$ws.Range("D38").Value = "'0.0164"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.10%  "
$ws.Range("E39").Value = "  -0.12%  "
$ws.Range("D40").Value = "'0.808"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.04%  "
$ws.Range("D41").Value = "'1.00"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.33%  "
$ws.Range("E42").Value = "  -0.68%  "
$ws.Range("E43").Value = "  +0.29%  "
$ws.Range("E44").Value = "  +3.39%  "
$ws.Range("D45").Value = "'63.57"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.35%  "
$ws.Range("D46").Value = "'1.74"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.52%  "
$ws.Range("D47").Value = "1.686.11"
$ws.Range("E47").Value = "  +0.35%  "
$ws.Range("E48").Value = "  +0.43%  "
$ws.Range("E49").Value = "  -0.06%  "
$ws.Range("E50").Value = "  +1.33%  "
